# "Set Exis Unit to 1 and MaxlineLoad 100%"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Storage")

# ExisUnits (column E, rows 7 and 10) -> 0
$ws.Range("E7").Value = 0
$ws.Range("E10").Value = 0

# MaxInvest / MaxlineLoad (column S, rows 7-11) -> 100%
$ws.Range("S7:S11").Value = 100
